# Daily auto push update: insert a new sample row for 2026/02/04 at the
# top of the "end of year" data block (row 755), pushing every following
# row down by one. The last existing row (796) is re-used by the shift,
# and a brand-new row 797 is created to hold the data that used to be
# the last row's trailing values (2027/01/05, 7:00).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 755..796 down to 756..797, making room for the new row 755.
$ws.Rows.Item(755).Insert()

# Column A holds a literal text date (e.g. "2026/12/29"), not a real
# Excel date value, in this workbook. Force text with a leading
# apostrophe so Excel doesn't reinterpret it as a date serial number,
# then strip the "number stored as text" formatting it applies so the
# new cell ends up with the same (default) style as its neighbours.
$ws.Range("A755").Value = "'2026/02/04"
$ws.Range("A755").ClearFormats()

$ws.Range("B755").Value = "水"
$ws.Range("C755").Value = 20
$ws.Range("D755").Value = 36
